$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original A3 style (style 2) in a scratch cell before we
# overwrite A3's style below, so we can re-apply it to the new row 4 cells.
$ws.Range("A3").Copy()
$ws.Range("Z100").PasteSpecial(-4122)

# Change A3's style to match A2 / header style (s=1), keeping its value.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Add new row 4 data
$ws.Range("A4").Value = "نامعتبر"
$ws.Range("B4").Value = "نامعتبریان"
$ws.Range("C4").Value = 123456

# Apply the preserved original "style 2" formatting to the new name cells.
$ws.Range("Z100").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)

# Clean up the scratch cell.
$ws.Range("Z100").Clear()
